$wb = $excel.ActiveWorkbook

# Corregir el "maximum change" usado en el calculo iterativo (calcPr/iterateDelta)
$excel.MaxChange = 0.0001

$newValues = @{
    8  = "NEWCALGUI"
    9  = "NEWCCLARI"
    10 = "NEWCGTAG1K"
    11 = "NEWCPTAG200G"
    12 = "NEW023400020002"
}

foreach ($ws in $wb.Worksheets) {
    foreach ($row in $newValues.Keys) {
        $ws.Cells.Item($row, 1).Value = $newValues[$row]
    }
}

$premium = $wb.Worksheets.Item("PREMIUM")
$premium.Activate()
$premium.Application.ActiveWindow.Zoom = 100
$premium.Range("C12").Select()

$mejorar = $wb.Worksheets.Item("MEJORAR")
$mejorar.Activate()
$mejorar.Application.ActiveWindow.Zoom = 85
$mejorar.Range("A11").Select()

$wb.Save()
